$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell G3: shared string "nur 25 funktioniert?"
$ws.Range("G3").Value = "nur 25 funktioniert?"

# Update data rows 34-57 (existing rows get new C/D/E values, A/B reshuffled)
# and add new data rows 58-89 (shifted/expanded data, new block A=5 and continuation blocks A=2,3,4)
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = -1
$ws.Range("D34").Value = 0.1
$ws.Range("E34").Value = 0.02
$ws.Range("A35").Value = 1
$ws.Range("B35").Value = 2
$ws.Range("C35").Value = -1
$ws.Range("D35").Value = 0.1
$ws.Range("E35").Value = 0.03
$ws.Range("A36").Value = 1
$ws.Range("B36").Value = 3
$ws.Range("C36").Value = -1
$ws.Range("D36").Value = 0.1
$ws.Range("E36").Value = 0.04
$ws.Range("A37").Value = 1
$ws.Range("B37").Value = 4
$ws.Range("C37").Value = -1
$ws.Range("D37").Value = 0.1
$ws.Range("E37").Value = 0.05
$ws.Range("A38").Value = 1
$ws.Range("B38").Value = 5
$ws.Range("C38").Value = -1
$ws.Range("D38").Value = 0.1
$ws.Range("E38").Value = 0.06
$ws.Range("A39").Value = 1
$ws.Range("B39").Value = 6
$ws.Range("C39").Value = -1
$ws.Range("D39").Value = 0.1
$ws.Range("E39").Value = 0.07
$ws.Range("A40").Value = 2
$ws.Range("B40").Value = 7
$ws.Range("C40").Value = -1
$ws.Range("D40").Value = 0.1
$ws.Range("E40").Value = 0.08
$ws.Range("A41").Value = 2
$ws.Range("B41").Value = 8
$ws.Range("C41").Value = -1
$ws.Range("D41").Value = 0.1
$ws.Range("E41").Value = 0.09
$ws.Range("A42").Value = 2
$ws.Range("B42").Value = 9
$ws.Range("C42").Value = -1
$ws.Range("D42").Value = 0.11
$ws.Range("E42").Value = 0.025
$ws.Range("A43").Value = 2
$ws.Range("B43").Value = 10
$ws.Range("C43").Value = -1
$ws.Range("D43").Value = 0.11
$ws.Range("E43").Value = 0.035
$ws.Range("A44").Value = 2
$ws.Range("B44").Value = 11
$ws.Range("C44").Value = -1
$ws.Range("D44").Value = 0.11
$ws.Range("E44").Value = 0.045
$ws.Range("A45").Value = 2
$ws.Range("B45").Value = 12
$ws.Range("C45").Value = -1
$ws.Range("D45").Value = 0.11
$ws.Range("E45").Value = 0.055
$ws.Range("A46").Value = 3
$ws.Range("B46").Value = 13
$ws.Range("C46").Value = -1
$ws.Range("D46").Value = 0.11
$ws.Range("E46").Value = 0.065
$ws.Range("A47").Value = 3
$ws.Range("B47").Value = 14
$ws.Range("C47").Value = -1
$ws.Range("D47").Value = 0.11
$ws.Range("E47").Value = 0.075
$ws.Range("A48").Value = 3
$ws.Range("B48").Value = 15
$ws.Range("C48").Value = -1
$ws.Range("D48").Value = 0.11
$ws.Range("E48").Value = 0.085
$ws.Range("A49").Value = 3
$ws.Range("B49").Value = 16
$ws.Range("C49").Value = -1
$ws.Range("D49").Value = 0.11
$ws.Range("E49").Value = 0.095
$ws.Range("A50").Value = 3
$ws.Range("B50").Value = 17
$ws.Range("C50").Value = -1
$ws.Range("D50").Value = 0.12
$ws.Range("E50").Value = 0.025
$ws.Range("A51").Value = 3
$ws.Range("B51").Value = 18
$ws.Range("C51").Value = -1
$ws.Range("D51").Value = 0.12
$ws.Range("E51").Value = 0.035
$ws.Range("A52").Value = 4
$ws.Range("B52").Value = 19
$ws.Range("C52").Value = -1
$ws.Range("D52").Value = 0.12
$ws.Range("E52").Value = 0.045
$ws.Range("A53").Value = 4
$ws.Range("B53").Value = 20
$ws.Range("C53").Value = -1
$ws.Range("D53").Value = 0.12
$ws.Range("E53").Value = 0.055
$ws.Range("A54").Value = 4
$ws.Range("B54").Value = 21
$ws.Range("C54").Value = -1
$ws.Range("D54").Value = 0.12
$ws.Range("E54").Value = 0.065
$ws.Range("A55").Value = 4
$ws.Range("B55").Value = 22
$ws.Range("C55").Value = -1
$ws.Range("D55").Value = 0.12
$ws.Range("E55").Value = 0.075
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = 23
$ws.Range("C56").Value = -1
$ws.Range("D56").Value = 0.12
$ws.Range("E56").Value = 0.085
$ws.Range("A57").Value = 4
$ws.Range("B57").Value = 24
$ws.Range("C57").Value = -1
$ws.Range("D57").Value = 0.12
$ws.Range("E57").Value = 0.095
$ws.Range("A58").Value = 5
$ws.Range("B58").Value = 1
$ws.Range("C58").Value = -3
$ws.Range("A59").Value = 5
$ws.Range("B59").Value = 2
$ws.Range("C59").Value = -3
$ws.Range("A60").Value = 5
$ws.Range("B60").Value = 3
$ws.Range("C60").Value = -3
$ws.Range("A61").Value = 5
$ws.Range("B61").Value = 4
$ws.Range("C61").Value = -3
$ws.Range("A62").Value = 5
$ws.Range("B62").Value = 5
$ws.Range("C62").Value = -3
$ws.Range("A63").Value = 5
$ws.Range("B63").Value = 6
$ws.Range("C63").Value = -3
$ws.Range("A64").Value = 5
$ws.Range("B64").Value = 7
$ws.Range("C64").Value = -3
$ws.Range("A65").Value = 5
$ws.Range("B65").Value = 8
$ws.Range("C65").Value = -3
$ws.Range("A66").Value = 2
$ws.Range("B66").Value = 9
$ws.Range("C66").Value = -3
$ws.Range("A67").Value = 2
$ws.Range("B67").Value = 10
$ws.Range("C67").Value = -3
$ws.Range("A68").Value = 2
$ws.Range("B68").Value = 11
$ws.Range("C68").Value = -3
$ws.Range("A69").Value = 2
$ws.Range("B69").Value = 12
$ws.Range("C69").Value = -3
$ws.Range("A70").Value = 2
$ws.Range("B70").Value = 13
$ws.Range("C70").Value = -3
$ws.Range("A71").Value = 2
$ws.Range("B71").Value = 14
$ws.Range("C71").Value = -3
$ws.Range("A72").Value = 2
$ws.Range("B72").Value = 15
$ws.Range("C72").Value = -3
$ws.Range("A73").Value = 2
$ws.Range("B73").Value = 16
$ws.Range("C73").Value = -3
$ws.Range("A74").Value = 3
$ws.Range("B74").Value = 17
$ws.Range("C74").Value = -3
$ws.Range("A75").Value = 3
$ws.Range("B75").Value = 18
$ws.Range("C75").Value = -3
$ws.Range("A76").Value = 3
$ws.Range("B76").Value = 19
$ws.Range("C76").Value = -3
$ws.Range("A77").Value = 3
$ws.Range("B77").Value = 20
$ws.Range("C77").Value = -3
$ws.Range("A78").Value = 3
$ws.Range("B78").Value = 21
$ws.Range("C78").Value = -3
$ws.Range("A79").Value = 3
$ws.Range("B79").Value = 22
$ws.Range("C79").Value = -3
$ws.Range("A80").Value = 3
$ws.Range("B80").Value = 23
$ws.Range("C80").Value = -3
$ws.Range("A81").Value = 3
$ws.Range("B81").Value = 24
$ws.Range("C81").Value = -3
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = 25
$ws.Range("C82").Value = -3
$ws.Range("A83").Value = 4
$ws.Range("B83").Value = 26
$ws.Range("C83").Value = -3
$ws.Range("A84").Value = 4
$ws.Range("B84").Value = 27
$ws.Range("C84").Value = -3
$ws.Range("A85").Value = 4
$ws.Range("B85").Value = 28
$ws.Range("C85").Value = -3
$ws.Range("A86").Value = 4
$ws.Range("B86").Value = 29
$ws.Range("C86").Value = -3
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = 30
$ws.Range("C87").Value = -3
$ws.Range("A88").Value = 4
$ws.Range("B88").Value = 31
$ws.Range("C88").Value = -3
$ws.Range("A89").Value = 4
$ws.Range("B89").Value = 32
$ws.Range("C89").Value = -3

# Update selection to match the authored edit (G57)
$ws.Range("G57").Select()
